# Lisää ELY Pääkäyttäjälle oikeus kaikkiin urakoihin
# Sähköposti Anne L Wednesday 31 August 2016 at 15:07 mukaisesti
# lisätään ELY Pääkäyttäjille kirj.oikeus kaikkiin urakoihin.
#
# Column J on the "Oikeudet" sheet is "ELY pääkäyttäjä". For every row in
# the "Urakat" (contracts) section its permission string is upgraded from
# "R*,W+" (write access limited to own organisation's contracts) to
# "R*,W*" (write access to ALL contracts). Row 24 ("Välitavoitteet") carries
# the extra "valmis" qualifier, so it goes from "R*,W+,valmis+" to
# "R*,W*,valmis*". Row 53 ("Ilmoitukset") gets the same base upgrade.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

$rows = @(7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,25,26,27,28,29,30,31,32,33,53)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Value = "R*,W*"
}

# Row 24 (Välitavoitteet) keeps its "valmis" suffix.
$ws.Cells.Item(24, 10).Value = "R*,W*,valmis*"

# Move the active selection on the frozen (bottom-right) pane up to J10,
# matching where the editor was actually working.
$ws.Range("J10").Select() | Out-Null
